$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2011 sbp_payment (B17) is recomputed with a flat 3x43.58 schedule instead
# of the previous half-weighted formula.
$ws.Range("B17").Formula = "=43.58+43.58+43.58"

# The highlighted/"classification_restrict" cell no longer needs the
# orange-tint theme fill -- reset it to plain white (Background 1).
$ws.Range("B17").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("B17").Interior.TintAndShade = 0

# Leave the cursor on the cell the author last looked at (E9) instead of B17.
$null = $ws.Range("E9").Select()
